# Apply the re-ordering of species observations across rows 2-10.
# The commit rotates each row's Id/Taxonsorteringsordning/.../Ost/Nord data
# so every row ends up holding the values that, before the edit, belonged
# to the 'next' record in the sequence (wrapping row 2 -> row 9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    [PSCustomObject]@{ Row = 2; A = 107856613; B = 96334; D = 'VU'; E = 220787; F = 'Knärot'; G = 'Goodyera repens'; H = '(L.) R. Br.'; I = $null; J = 'plantor/tuvor'; Q = 585871.3090144017; R = 6393989.072683465 }
    [PSCustomObject]@{ Row = 3; A = 107856240; B = 103265; D = 'LC'; E = 221144; F = 'Grönpyrola'; G = 'Pyrola chlorantha'; H = 'Sw.'; I = '10'; J = 'plantor/tuvor'; Q = 585839.5429383111; R = 6393945.609328114 }
    [PSCustomObject]@{ Row = 4; A = 107856215; B = 96334; D = 'VU'; E = 220787; F = 'Knärot'; G = 'Goodyera repens'; H = '(L.) R. Br.'; I = $null; J = $null; Q = 585827.0278735023; R = 6393929.297724779 }
    [PSCustomObject]@{ Row = 5; A = 107856294; B = 96334; D = 'VU'; E = 220787; F = 'Knärot'; G = 'Goodyera repens'; H = '(L.) R. Br.'; I = $null; J = 'plantor/tuvor'; Q = 585820.3320085882; R = 6393992.268917265 }
    [PSCustomObject]@{ Row = 6; A = 107856266; B = 96334; D = 'VU'; E = 220787; F = 'Knärot'; G = 'Goodyera repens'; H = '(L.) R. Br.'; I = '5'; J = 'plantor/tuvor'; Q = 585831.1129932627; R = 6393989.288691249 }
    [PSCustomObject]@{ Row = 7; A = 107856259; B = 96334; D = 'VU'; E = 220787; F = 'Knärot'; G = 'Goodyera repens'; H = '(L.) R. Br.'; I = $null; J = 'plantor/tuvor'; Q = 585838.8369974099; R = 6393953.617245103 }
    [PSCustomObject]@{ Row = 8; A = 107856280; B = 103265; D = 'LC'; E = 221144; F = 'Grönpyrola'; G = 'Pyrola chlorantha'; H = 'Sw.'; I = $null; J = 'plantor/tuvor'; Q = 585819.2153204344; R = 6393994.384627678 }
    [PSCustomObject]@{ Row = 9; A = 107856371; B = 96334; D = 'VU'; E = 220787; F = 'Knärot'; G = 'Goodyera repens'; H = '(L.) R. Br.'; I = $null; J = 'plantor/tuvor'; Q = 585757.3301987207; R = 6393981.30443168 }
    [PSCustomObject]@{ Row = 10; A = 107856305; B = 103265; D = 'LC'; E = 221144; F = 'Grönpyrola'; G = 'Pyrola chlorantha'; H = 'Sw.'; I = $null; J = $null; Q = 585775.987021131; R = 6393985.979127978 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A    # Id
    $ws.Cells.Item($r.Row, 2).Value = $r.B    # Taxonsorteringsordning
    $ws.Cells.Item($r.Row, 4).Value = $r.D    # Rödlistade
    $ws.Cells.Item($r.Row, 5).Value = $r.E    # TaxonId
    $ws.Cells.Item($r.Row, 6).Value = $r.F    # Artnamn
    $ws.Cells.Item($r.Row, 7).Value = $r.G    # Vetenskapligt namn
    $ws.Cells.Item($r.Row, 8).Value = $r.H    # Auktor

    # Antal (I) is numeric-looking text ('10', '5', or blank) - keep it text
    if ($null -eq $r.I) {
        $ws.Cells.Item($r.Row, 9).Value = ''
    } else {
        $ws.Cells.Item($r.Row, 9).Value = "'" + $r.I
    }

    # Enhet (J) is plain text, sometimes blank
    if ($null -eq $r.J) {
        $ws.Cells.Item($r.Row, 10).Value = ''
    } else {
        $ws.Cells.Item($r.Row, 10).Value = $r.J
    }

    $ws.Cells.Item($r.Row, 17).Value = $r.Q   # Ost
    $ws.Cells.Item($r.Row, 18).Value = $r.R   # Nord
}
